# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# Two pairs of match rows had been entered in the wrong order; this swaps
# the full record (every column except the leading running-number column A)
# between the two rows in each pair so the data lines up with the correct
# match id / odds set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB, $firstCol, $lastCol) {
    $valsA = @()
    $valsB = @()

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $valsA += , ($ws.Cells.Item($rowA, $c).Value2())
        $valsB += , ($ws.Cells.Item($rowB, $c).Value2())
    }

    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($rowA, $c).Value2 = $valsB[$i]
        $ws.Cells.Item($rowB, $c).Value2 = $valsA[$i]
        $i = $i + 1
    }
}

# Columns B (2) through AD (30); column A (running index) is left untouched.
Swap-Rows 136 137 2 30
Swap-Rows 255 256 2 30
